$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna2"
$ws.Range("C2").Value = "Epha2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.408030333333333
$ws.Range("H2").Value = 4.224091
$ws.Range("I2").Value = 0.3454737251382253
$ws.Range("J2").Value = 0.3454737251382253
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 11.55727433333333
$ws.Range("N2").Value = 34.671823
$ws.Range("O2").Value = 0.5239815261112396
$ws.Range("P2").Value = 0.5239815261112395
$ws.Range("Q2").Value = 16.27299283198811
$ws.Range("R2").Value = 146.456935487893
$ws.Range("S2").Value = 0.1810218497292622
$ws.Range("T2").Value = 0.1810218497292622

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna2"
$ws.Range("C3").Value = "Epha2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.408030333333333
$ws.Range("H3").Value = 4.224091
$ws.Range("I3").Value = 0.3454737251382253
$ws.Range("J3").Value = 0.3454737251382253
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.09477133333333332
$ws.Range("N3").Value = 0.284314
$ws.Range("O3").Value = 0.004296724853919303
$ws.Range("P3").Value = 0.004296724853919302
$ws.Range("Q3").Value = 0.1334409120637778
$ws.Range("R3").Value = 1.200968208574
$ws.Range("S3").Value = 0.001484405541177498
$ws.Range("T3").Value = 0.001484405541177498

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna2"
$ws.Range("C4").Value = "Epha2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.408030333333333
$ws.Range("H4").Value = 4.224091
$ws.Range("I4").Value = 0.3454737251382253
$ws.Range("J4").Value = 0.3454737251382253
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.404599
$ws.Range("N4").Value = 31.213797
$ws.Range("O4").Value = 0.471721749034841
$ws.Range("P4").Value = 0.471721749034841
$ws.Range("Q4").Value = 14.64999099816966
$ws.Range("R4").Value = 131.849918983527
$ws.Range("S4").Value = 0.1629674698677856
$ws.Range("T4").Value = 0.1629674698677856

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna2"
$ws.Range("C5").Value = "Epha2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.015377
$ws.Range("H5").Value = 6.046131
$ws.Range("I5").Value = 0.494492045565236
$ws.Range("J5").Value = 0.4944920455652361
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.55727433333333
$ws.Range("N5").Value = 34.671823
$ws.Range("O5").Value = 0.5239815261112396
$ws.Range("P5").Value = 0.5239815261112395
$ws.Range("Q5").Value = 23.29226487409034
$ws.Range("R5").Value = 209.630383866813
$ws.Range("S5").Value = 0.259104696685141
$ws.Range("T5").Value = 0.259104696685141

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna2"
$ws.Range("C6").Value = "Epha2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.015377
$ws.Range("H6").Value = 6.046131
$ws.Range("I6").Value = 0.494492045565236
$ws.Range("J6").Value = 0.4944920455652361
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09477133333333332
$ws.Range("N6").Value = 0.284314
$ws.Range("O6").Value = 0.004296724853919303
$ws.Range("P6").Value = 0.004296724853919302
$ws.Range("Q6").Value = 0.1909999654593333
$ws.Range("R6").Value = 1.718999689134
$ws.Range("S6").Value = 0.002124696262245546
$ws.Range("T6").Value = 0.002124696262245546

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna2"
$ws.Range("C7").Value = "Epha2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.015377
$ws.Range("H7").Value = 6.046131
$ws.Range("I7").Value = 0.494492045565236
$ws.Range("J7").Value = 0.4944920455652361
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.404599
$ws.Range("N7").Value = 31.213797
$ws.Range("O7").Value = 0.471721749034841
$ws.Range("P7").Value = 0.471721749034841
$ws.Range("Q7").Value = 20.969189518823
$ws.Range("R7").Value = 188.722705669407
$ws.Range("S7").Value = 0.2332626526178495
$ws.Range("T7").Value = 0.2332626526178495

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna2"
$ws.Range("C8").Value = "Epha2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6522436666666667
$ws.Range("H8").Value = 1.956731
$ws.Range("I8").Value = 0.1600342292965385
$ws.Range("J8").Value = 0.1600342292965385
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.55727433333333
$ws.Range("N8").Value = 34.671823
$ws.Range("O8").Value = 0.5239815261112396
$ws.Range("P8").Value = 0.5239815261112395
$ws.Range("Q8").Value = 7.53815898784589
$ws.Range("R8").Value = 67.84343089061301
$ws.Range("S8").Value = 0.08385497969683632
$ws.Range("T8").Value = 0.08385497969683629

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna2"
$ws.Range("C9").Value = "Epha2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6522436666666667
$ws.Range("H9").Value = 1.956731
$ws.Range("I9").Value = 0.1600342292965385
$ws.Range("J9").Value = 0.1600342292965385
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09477133333333332
$ws.Range("N9").Value = 0.284314
$ws.Range("O9").Value = 0.004296724853919303
$ws.Range("P9").Value = 0.004296724853919302
$ws.Range("Q9").Value = 0.06181400194822221
$ws.Range("R9").Value = 0.5563260175339999
$ws.Range("S9").Value = 0.0006876230504962578
$ws.Range("T9").Value = 0.0006876230504962577

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna2"
$ws.Range("C10").Value = "Epha2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6522436666666667
$ws.Range("H10").Value = 1.956731
$ws.Range("I10").Value = 0.1600342292965385
$ws.Range("J10").Value = 0.1600342292965385
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.404599
$ws.Range("N10").Value = 31.213797
$ws.Range("O10").Value = 0.471721749034841
$ws.Range("P10").Value = 0.471721749034841
$ws.Range("Q10").Value = 6.786333801956333
$ws.Range("R10").Value = 61.077004217607
$ws.Range("S10").Value = 0.07549162654920596
$ws.Range("T10").Value = 0.07549162654920596
